# The "Learning Resources" cell for the first DVS section lists several
# articles/books. One entry - a Harvard Business Review article - had its
# URL pasted in as several broken/hyphen-mangled fragments separated by
# manual line breaks:
#
#   ...You've Already Made by Kevin Troyanos
#   https://hbr.org/2018/10/how-
#   tomake-sure-youre-not-
#   usingdata-just-to-justify-
#   decisionsyouve-already-made
#
# This script removes that broken URL (everything after "Troyanos ")
# while leaving the author's name and the rest of the paragraph intact.
# Only the first occurrence of this citation block is edited (a second,
# identical copy of the same reference list elsewhere in the document is
# left untouched).

$d = $word.ActiveDocument

# Locate the end of "Troyanos" - the author's surname that introduces the
# broken URL we want to strip out.
$anchorStart = $d.Content
$findStart = $anchorStart.Find
$findStart.ClearFormatting()
$findStart.Text = "Troyanos"
$findStart.Forward = $true
$findStart.Wrap = 1  # wdFindStop
$foundStart = $findStart.Execute()

if (-not $foundStart) {
    throw "Could not find 'Troyanos' anchor text."
}

# Keep the single non-breaking space that immediately follows the name;
# deletion begins at the manual line break right after it.
$deleteStart = $anchorStart.End + 1

# Locate the tail end of the mangled URL text ("...-already-made").
$anchorEnd = $d.Content
$findEnd = $anchorEnd.Find
$findEnd.ClearFormatting()
$findEnd.Text = "decisionsyouve-already-made"
$findEnd.Forward = $true
$findEnd.Wrap = 1  # wdFindStop
$foundEnd = $findEnd.Execute()

if (-not $foundEnd) {
    throw "Could not find end of mangled URL text."
}

# Include the trailing non-breaking space (the paragraph's "end of
# paragraph" marker run) so nothing but the closing paragraph mark is left.
$deleteEnd = $anchorEnd.End + 1

$deadLinkRange = $d.Range($deleteStart, $deleteEnd)
$deadLinkRange.Delete()

Write-Output "Removed mangled hbr.org URL text ($deleteStart-$deleteEnd)."
